$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "rowname"
$ws.Range("B1").Value = "estimate"
$ws.Range("C1").Value = "SE"
$ws.Range("D1").Value = "df"
$ws.Range("E1").Value = "t.ratio"
$ws.Range("F1").Value = "p.value"
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"
$ws.Range("A2").Value = "r_xy1y2"
$ws.Range("B2").Value = -0.533680282822207
$ws.Range("C2").Value = 0.106303989026289
$ws.Range("D2").Value = 64.0692747366915
$ws.Range("E2").Value = -5.02032226363801
$ws.Range("F2").Value = 0.00000437462177627036
$ws.Range("G2").Value = -0.746042491862588
$ws.Range("H2").Value = -0.321318073781826
$ws.Range("A3").Value = "w_11"
$ws.Range("B3").Value = -0.119972253958424
$ws.Range("C3").Value = 0.0412243661702373
$ws.Range("D3").Value = 68.8824312298728
$ws.Range("E3").Value = -2.91022676887244
$ws.Range("F3").Value = 0.00486049239182246
$ws.Range("G3").Value = -0.202215118496997
$ws.Range("H3").Value = -0.0377293894198507
$ws.Range("A4").Value = "w_21"
$ws.Range("B4").Value = -0.0562886813695994
$ws.Range("C4").Value = 0.0402521284260808
$ws.Range("D4").Value = 69.0816253252937
$ws.Range("E4").Value = -1.39840260802527
$ws.Range("F4").Value = 0.166465416443895
$ws.Range("G4").Value = -0.136587787757276
$ws.Range("H4").Value = 0.0240104250180769
$ws.Range("A5").Value = "r_xy1"
$ws.Range("B5").Value = -0.331034594329411
$ws.Range("C5").Value = 0.113748728404993
$ws.Range("D5").Value = 68.8824312298728
$ws.Range("E5").Value = -2.91022676887244
$ws.Range("F5").Value = 0.00486049239182247
$ws.Range("G5").Value = -0.557964008429192
$ws.Range("H5").Value = -0.104105180229629
$ws.Range("A6").Value = "r_xy2"
$ws.Range("B6").Value = -0.165973017005494
$ws.Range("C6").Value = 0.118687576848752
$ws.Range("D6").Value = 69.0816253252937
$ws.Range("E6").Value = -1.39840260802527
$ws.Range("F6").Value = 0.166465416443895
$ws.Range("G6").Value = -0.402743263273968
$ws.Range("H6").Value = 0.0707972292629808
$ws.Range("A7").Value = "b_11"
$ws.Range("B7").Value = -0.342015842091268
$ws.Range("C7").Value = 0.117522059019401
$ws.Range("D7").Value = 68.8824312298728
$ws.Range("E7").Value = -2.91022676887244
$ws.Range("F7").Value = 0.00486049239182247
$ws.Range("G7").Value = -0.57647307401846
$ws.Range("H7").Value = -0.107558610164075
$ws.Range("A8").Value = "b_21"
$ws.Range("B8").Value = -0.160467275754461
$ws.Range("C8").Value = 0.114750412244341
$ws.Range("D8").Value = 69.0816253252937
$ws.Range("E8").Value = -1.39840260802527
$ws.Range("F8").Value = 0.166465416443895
$ws.Range("G8").Value = -0.389383259110703
$ws.Range("H8").Value = 0.0684487076017813
$ws.Range("A9").Value = "main_effect"
$ws.Range("B9").Value = -0.0881304676640117
$ws.Range("C9").Value = 0.040244413362567
$ws.Range("D9").Value = 69.040710027208
$ws.Range("E9").Value = -2.18988078842231
$ws.Range("F9").Value = 0.0319143445664243
$ws.Range("G9").Value = -0.168415031034521
$ws.Range("H9").Value = -0.00784590429350256
$ws.Range("A10").Value = "moderator_effect"
$ws.Range("B10").Value = -0.104209736816482
$ws.Range("C10").Value = 0.012766735176762
$ws.Range("D10").Value = 67.9241855539697
$ws.Range("E10").Value = -8.16259876731554
$ws.Range("F10").Value = 0.0000000000111846009858307
$ws.Range("G10").Value = -0.129685875097318
$ws.Range("H10").Value = -0.078733598535645
$ws.Range("A11").Value = "interaction"
$ws.Range("B11").Value = -0.0636835725888246
$ws.Range("C11").Value = 0.0126851562996428
$ws.Range("D11").Value = 64.0692747366915
$ws.Range("E11").Value = -5.02032226363801
$ws.Range("F11").Value = 0.00000437462177627037
$ws.Range("G11").Value = -0.0890245578000989
$ws.Range("H11").Value = -0.0383425873775503
$ws.Range("A12").Value = "q_b11_b21"
$ws.Range("B12").Value = -0.194507348435464
$ws.Range("A13").Value = "q_rxy1_rxy2"
$ws.Range("B13").Value = -0.176466995706537
$ws.Range("A14").Value = "cross_over_point"
$ws.Range("B14").Value = -1.63636763108935
$ws.Range("A15").Value = "interaction_vs_main"
$ws.Range("B15").Value = -0.0244468950751872
$ws.Range("C15").Value = 0.0412469571475937
$ws.Range("D15").Value = 69.0223993179759
$ws.Range("E15").Value = -0.592695722685894
$ws.Range("F15").Value = 0.555321724222152
$ws.Range("G15").Value = -0.106731846733835
$ws.Range("H15").Value = 0.0578380565834612
$ws.Range("A16").Value = "interaction_vs_main_bscale"
$ws.Range("B16").Value = -0.0696929925860571
$ws.Range("C16").Value = 0.117586461178144
$ws.Range("D16").Value = 69.0223993179759
$ws.Range("E16").Value = -0.592695722685895
$ws.Range("F16").Value = 0.555321724222152
$ws.Range("G16").Value = -0.304270206103481
$ws.Range("H16").Value = 0.164884220931367
$ws.Range("A17").Value = "interaction_vs_main_rscale"
$ws.Range("B17").Value = -0.0834422283435352
$ws.Range("C17").Value = 0.125120497487471
$ws.Range("D17").Value = 69.0319993883393
$ws.Range("E17").Value = -0.66689495341793
$ws.Range("F17").Value = 0.507062244495859
$ws.Range("G17").Value = -0.33304872697581
$ws.Range("H17").Value = 0.16616427028874
$ws.Range("A18").Value = "dadas"
$ws.Range("B18").Value = -0.112577362739199
$ws.Range("C18").Value = 0.0805042568521617
$ws.Range("D18").Value = 69.0816253252937
$ws.Range("E18").Value = -1.39840260802527
$ws.Range("F18").Value = 0.916767291778053
$ws.Range("G18").Value = -0.273175575514552
$ws.Range("H18").Value = 0.0480208500361539
$ws.Range("A19").Value = "dadas_bscale"
$ws.Range("B19").Value = -0.320934551508921
$ws.Range("C19").Value = 0.229500824488681
$ws.Range("D19").Value = 69.0816253252937
$ws.Range("E19").Value = -1.39840260802527
$ws.Range("F19").Value = 0.916767291778053
$ws.Range("G19").Value = -0.778766518221405
$ws.Range("H19").Value = 0.136897415203563
$ws.Range("A20").Value = "dadas_rscale"
$ws.Range("B20").Value = -0.331946034010987
$ws.Range("C20").Value = 0.237375153697503
$ws.Range("D20").Value = 69.0816253252937
$ws.Range("E20").Value = -1.39840260802527
$ws.Range("F20").Value = 0.916767291778053
$ws.Range("G20").Value = -0.805486526547936
$ws.Range("H20").Value = 0.141594458525962
$ws.Range("A21").Value = "abs_diff"
$ws.Range("B21").Value = 0.0636835725888246
$ws.Range("C21").Value = 0.0126851562996428
$ws.Range("D21").Value = 64.0692747366915
$ws.Range("E21").Value = 5.02032226363801
$ws.Range("F21").Value = 0.00000218731088813518
$ws.Range("G21").Value = 0.0383425873775503
$ws.Range("H21").Value = 0.0890245578000989
$ws.Range("A22").Value = "abs_sum"
$ws.Range("B22").Value = 0.176260935328023
$ws.Range("C22").Value = 0.0804888267251341
$ws.Range("D22").Value = 69.040710027208
$ws.Range("E22").Value = 2.18988078842231
$ws.Range("F22").Value = 0.0159571722832122
$ws.Range("G22").Value = 0.0156918085870051
$ws.Range("H22").Value = 0.336830062069042
$ws.Range("A23").Value = "abs_diff_bscale"
$ws.Range("B23").Value = 0.181548566336807
$ws.Range("C23").Value = 0.0361627315544573
$ws.Range("D23").Value = 64.0692747366915
$ws.Range("E23").Value = 5.02032226363801
$ws.Range("F23").Value = 0.00000218731088813518
$ws.Range("G23").Value = 0.10930670948664
$ws.Range("H23").Value = 0.253790423186974
$ws.Range("A24").Value = "abs_sum_bscale"
$ws.Range("B24").Value = 0.502483117845728
$ws.Range("C24").Value = 0.229456836418817
$ws.Range("D24").Value = 69.040710027208
$ws.Range("E24").Value = 2.18988078842231
$ws.Range("F24").Value = 0.0159571722832122
$ws.Range("G24").Value = 0.0447340693430617
$ws.Range("H24").Value = 0.960232166348395
$ws.Range("A25").Value = "abs_diff_rscale"
$ws.Range("B25").Value = 0.165061577323917
$ws.Range("C25").Value = 0.0364127111697027
$ws.Range("D25").Value = 64.6270860321452
$ws.Range("E25").Value = 4.53307573156642
$ws.Range("F25").Value = 0.0000128756838171037
$ws.Range("G25").Value = 0.0923324157042769
$ws.Range("H25").Value = 0.237790738943557
$ws.Range("A26").Value = "abs_sum_rscale"
$ws.Range("B26").Value = 0.497007611334904
$ws.Range("C26").Value = 0.229619560772235
$ws.Range("D26").Value = 69.0438268841139
$ws.Range("E26").Value = 2.16448289363247
$ws.Range("F26").Value = 0.0169439365394003
$ws.Range("G26").Value = 0.0389343087365633
$ws.Range("H26").Value = 0.955080913933245

# Create empty (but present) cells to match the exported rectangular grid
$ws.Range("G12").Font.Bold = $false
$ws.Range("H12").Font.Bold = $false
$ws.Range("G13").Font.Bold = $false
$ws.Range("H13").Font.Bold = $false
$ws.Range("G14").Font.Bold = $false
$ws.Range("H14").Font.Bold = $false
